# Grupo 1 - Notas: build the grade sheet (TP1-TP4 comments/grades for
# "ANDREW E GABRIEL AMORIM") into the previously-empty Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Values / formulas
# ---------------------------------------------------------------------

$ws.Range("B1").Value = "ANDREW E GABRIEL AMORIM"

$ws.Range("B2").Value = "TOTAL"
$ws.Range("C2").Formula = "=SUM(C4:C999)"

$ws.Range("B4").Value = "Comentários"
$ws.Range("C4").Value = "Nota"

$longComment = "AEC para herança e nomenclatura (andrew e gabriel apresentou)`n- boa contextualização`n- boa explicação de ambuiguidade`n- boa explicação de diamante`n- padrões de nomenclatura ok.`n- Typo: CAPITALIZED`n- AntLR4 + Python`n- Excelente ambos apresentarem, mas o Gabriel tinha que assistir antes o do Andrew, pois repetiu muito.`n- microfone Gabriel ruim. Não entendi algumas palavras. Muito ruído.`n- código tem *várias* oportunidades de extract method, não? Bem `"linguição`".`n- Falta um fechamento da apresentação. Final abrupto."

$ws.Range("A5").Value = "TP1"
$ws.Range("B5").Value = $longComment
$ws.Range("C5").Value = 10

$ws.Range("A6").Value = "TP2"

$ws.Range("A7").Value = "TP3"

$ws.Range("A8").Value = "TP4"

Write-Output "values written"
